$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown Chart")

# Row 6: I6 goes from text "0.5" shared-string to numeric 0.5; K6 gets numeric 0.5
$ws.Range("I6").Value = 0.5
$ws.Range("K6").Value = 0.5

# Row 7: H7, I7, J7 go from text "0.5" to numeric 0.5; K7 gets numeric 0.5
$ws.Range("H7").Value = 0.5
$ws.Range("I7").Value = 0.5
$ws.Range("J7").Value = 0.5
$ws.Range("K7").Value = 0.5

# Update the active selection to L12 to match the saved cursor position
$ws.Range("L12").Select()
